$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 15: give it the "last row of group" border style (same as rows 4, 11) ---
# Copy formatting (styles) from row 4 (which already uses style 6/7, the bottom-border variant)
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Prepare rows 16 and 17 with the "normal" group style (same as rows 5, 7, 9, 12, 14) ---
$ws.Range("A5:E5").Copy() | Out-Null
$ws.Range("A16:E17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Fill in the new cell values (order matters so new shared strings are appended in the
#     same sequence as the target workbook) ---
$ws.Range("A16").Value = "SCRIPT/T01P02A/us0408.ssb"
$ws.Range("C16").Value = " There are exploration teams\neven more famous than Team [CS:X]Charm[CR]."
$ws.Range("C17").Value = " The most famous would have to\nbe Team [CS:X]Raider[CR]."
$ws.Range("D16").Value = " Есть более известная команда\nисследователей, даже больше чем Команда\n[CS:X]Шарм[CR]."
$ws.Range("D17").Value = " Самой известной командой\nявляется Команда [CS:X]Рейдер[CR]."
$ws.Range("E16").Value = " Åòóû áïìåå éèâåòóîàÿ ëïíàîäà\néòòìåäïâàóåìåê, äàçå áïìûšå œåí Ëïíàîäà\n[CS:X]Šàñí[CR]."
$ws.Range("E17").Value = " Òàíïê éèâåòóîïê ëïíàîäïê\nÿâìÿåóòÿ Ëïíàîäà [CS:X]Ñåêäåñ[CR]."
$ws.Range("A17").Value = "SCRIPT/T01P02A/us2004.ssb"

$ws.Range("B16").Value = 38
$ws.Range("B17").Value = 41

# --- Row heights to match the multi-line wrapped text (same as other similar rows) ---
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(17).RowHeight = 43.2

# --- Update selection to match the new last cell ---
$ws.Range("E17").Select()
